$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices for B,C,D,E,F,I,J,K,L,M,N (G and H are skipped)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @{}
$data[2] = @(1.02, 1.026334791081971, 1.028979306322998, 1.026571595137207, 1.036558711088824, 1.032847858321424, 1.031498808557979, 1.031794549883021, 1.029393852075565, 1.039352107833707, 1.032963655180705)
$data[3] = @(1.02, 1.027381719341492, 1.029939222441116, 1.027462349976322, 1.03781154444677, 1.033059745142316, 1.032184951014921, 1.032562105923905, 1.030091923992617, 1.040413392444599, 1.033650772038632)
$data[4] = @(1.02, 1.028058887736735, 1.030560379069004, 1.028038890305033, 1.038621865134681, 1.033194749697788, 1.032628117825693, 1.033058178294346, 1.030543155964259, 1.041099217436501, 1.034094568197106)
$data[5] = @(1.02, 1.028343506583088, 1.030821519656093, 1.028281306513553, 1.038962442441235, 1.033251002594396, 1.032814230483779, 1.033266586415766, 1.030732741951685, 1.041387324541319, 1.034280945156514)
$data[6] = @(1.02, 1.028391291674742, 1.030865366679844, 1.028322011521795, 1.039019622129945, 1.033260418200577, 1.032845468203012, 1.033301570805248, 1.030764567707025, 1.041435686489544, 1.034312227236892)
$data[7] = @(1.02, 1.028062691073626, 1.03056386841748, 1.028042129330685, 1.038626416263893, 1.03319550332851, 1.032630605437048, 1.033060963606755, 1.030545689661436, 1.041103067976234, 1.034097059341154)
$data[8] = @(1.02, 1.026688660396452, 1.02930370890895, 1.026872596152615, 1.036982184522951, 1.032919901256673, 1.031730862016554, 1.032054070731334, 1.029629865449779, 1.039710960326064, 1.033196038181801)
$data[9] = @(1.02, 1.024265396061721, 1.027083352812978, 1.024812981816964, 1.03408212312184, 1.032418176092984, 1.030139168811194, 1.030275294860715, 1.028012490217303, 1.037250969979556, 1.03160208458952)
$data[10] = @(1.02, 1.022648474331617, 1.025603252066035, 1.023440756536765, 1.032146831783656, 1.032072889908348, 1.029073845297226, 1.029086412953975, 1.026931836242771, 1.035606253999995, 1.030535248193968)
$data[11] = @(1.02, 1.021947985318461, 1.024962382946992, 1.022846770518303, 1.031308353060351, 1.031920815943293, 1.028611551268247, 1.028570892949972, 1.02646333013844, 1.034892936918904, 1.030072297654442)
$data[12] = @(1.02, 1.021687738833578, 1.024724339009861, 1.022626167155337, 1.030996829911168, 1.031863943866995, 1.028439683849827, 1.028379296459797, 1.026289219424889, 1.034627805710474, 1.029900186164574)
$data[13] = @(1.02, 1.021743565031621, 1.02477540009335, 1.022673485991779, 1.031063656101791, 1.031876160544247, 1.02847655682498, 1.028420399512983, 1.026326570689661, 1.034684685130703, 1.029937111503582)
$data[14] = @(1.02, 1.021926474363087, 1.024942706093296, 1.022828534774165, 1.031282603997434, 1.031916122734385, 1.028597347732561, 1.028555057757907, 1.026448939861976, 1.034871024638711, 1.030058073948106)
$data[15] = @(1.02, 1.022039163699103, 1.025045789282771, 1.022924069352728, 1.031417495048514, 1.031940693737472, 1.028671751001279, 1.028638010640558, 1.026524324055503, 1.034985811523792, 1.030132582878)
$data[16] = @(1.02, 1.022694956003058, 1.025645784928368, 1.023480181570888, 1.032202468471935, 1.032082928561507, 1.029104505094259, 1.029120610992119, 1.026962917303392, 1.035653570263697, 1.030565951531429)
$data[17] = @(1.02, 1.02310622251664, 1.026022152664404, 1.023829068676737, 1.032694730371736, 1.032171462580168, 1.029375691849058, 1.029423138791745, 1.027237880978557, 1.036072130231454, 1.030837523402526)
$data[18] = @(1.02, 1.023346073488265, 1.026241683956924, 1.024032587593704, 1.032981811766679, 1.032222855706323, 1.029533773845519, 1.029599528082124, 1.027398207094141, 1.036316158886284, 1.030995829893554)
$data[19] = @(1.02, 1.023427850760998, 1.026316538862431, 1.024101985503935, 1.033079691252178, 1.032240337493217, 1.029587659336841, 1.02965966039767, 1.027452864740412, 1.036399347691601, 1.031049791908455)
$data[20] = @(1.02, 1.02306210101721, 1.025981771725849, 1.023791634430033, 1.032641920167072, 1.032161989293999, 1.029346606075717, 1.029390687670173, 1.02720838572478, 1.036027234127349, 1.030808396324051)
$data[21] = @(1.02, 1.021872613612912, 1.024893438554887, 1.022782875958518, 1.031218131380623, 1.031904365494358, 1.028561781987122, 1.028515407283589, 1.026412907585135, 1.034816157099747, 1.030022457695229)
$data[22] = @(1.02, 1.021124424663151, 1.024209179764084, 1.02214880007122, 1.030322504995983, 1.031740159485851, 1.028067459105764, 1.027964450618846, 1.025912257272852, 1.034053700388877, 1.02952743281868)
$data[23] = @(1.02, 1.021521083566645, 1.024571916538945, 1.022484919541873, 1.030797335278973, 1.031827419362711, 1.028329591856397, 1.028256583285378, 1.026177709086794, 1.034457988864861, 1.029789937827887)
$data[24] = @(1.02, 1.0230820377278, 1.026000018126308, 1.023808549296011, 1.032665782968899, 1.032166270628167, 1.029359748983143, 1.029405351151433, 1.027221713531811, 1.036047521084544, 1.030821557895912)
$data[25] = @(1.02, 1.02489211420837, 1.027657343827511, 1.025345292115247, 1.034832189641107, 1.032549788680284, 1.030551398544891, 1.030735684989623, 1.02843104417322, 1.037887763587086, 1.032014899736722)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value2 = $vals[$i]
    }
}
